# Insert a new record row at row 85 (pushing existing rows 85..99 down to 86..100)
# and populate it with the new weekly data point for Jengibre.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 85; this shifts rows 85-99 down to 86-100.
$ws.Rows.Item(85).Insert()

# Fill the new row 85 with the new record.
$ws.Cells.Item(85, 1).Value = 8
$ws.Cells.Item(85, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(85, 3).Value = "Coquimbo"
$ws.Cells.Item(85, 4).Value = 44943
$ws.Cells.Item(85, 4).NumberFormat = $ws.Cells.Item(86, 4).NumberFormat
$ws.Cells.Item(85, 5).Value = 4
$ws.Cells.Item(85, 6).Value = 100114007
$ws.Cells.Item(85, 7).Value = "Jengibre"
$ws.Cells.Item(85, 8).Value = "Sin especificar"
$ws.Cells.Item(85, 9).Value = "Primera"
$ws.Cells.Item(85, 10).Value = 400
$ws.Cells.Item(85, 11).Value = 14000
$ws.Cells.Item(85, 12).Value = 15000
$ws.Cells.Item(85, 13).Value = 14500
$ws.Cells.Item(85, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(85, 15).Value = "Perú"
$ws.Cells.Item(85, 16).Value = 1115
$ws.Cells.Item(85, 17).Value = 13
$ws.Cells.Item(85, 18).Value = "Hortaliza"
